# Datebook_jph.xlsx - "add 0705 comment, 0706 goal"
#
# Row 11 (date 44379 = 2021-07-02) : C11 time re-entered as text "08:27/18:03"
# Row 12 (date 44382 = 2021-07-05) : C12 time added, D12 goal gets a 3rd bullet,
#                                     E12 gets the long 0705 comment text
# Row 13 (date 44383 = 2021-07-06) : D13 gets the short 0706 goal text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 : normalize the time value to literal text (matches the rest of column C) ---
$ws.Range("C11").Value = "08:27`n18:03"

# --- Row 12 : fill in the 07/05 time, goal (3rd bullet) and comment ---
$ws.Range("C12").Value = "08:25`n18:02"

$ws.Range("D12").Value = "1. 인터페이스 구상도 작성`n2. 요구사항 분류하여 수정`n3. 추가 자료조사"

$ws.Range("E12").Value = "1. 요구사항 분석 (실현성 확인)`n2. 데이터 흐름도 수정`n - 차량 관련 새소식 창 삭제 및 다른 사항 수정`n3. Menu_Tree 작성`n4. Flow_Chart 작성`n5. 필요 자료 추가 조사`n -> 주유, 부품들의 수명정보(?) 어떻게 자동으로 받아올지 알아봐야함 -> 블루투스 => 블루투스 연동 어떻게?"

# --- Row 13 : fill in the 07/06 goal ---
$ws.Range("D13").Value = "추가 자료조사"

# --- formatting to match the rest of the sheet ---
# C11/C12 are time-ish text that wraps onto two lines, center aligned (like the
# other cells in column C), E12 wraps like the other Comment cells in column E.
$ws.Range("C11:C12").WrapText = $true
$ws.Range("C11:C12").HorizontalAlignment = -4108   # xlCenter
$ws.Range("C11:C12").VerticalAlignment = -4108     # xlCenter

$ws.Range("E12").WrapText = $true
$ws.Range("E12").HorizontalAlignment = -4131       # xlLeft
$ws.Range("E12").VerticalAlignment = -4108         # xlCenter

# Row 12 grows to fit the new, much longer comment/goal text.
$ws.Rows.Item(12).RowHeight = 105.4

# --- view state: last-edited cell was D12 ---
$ws.Range("D12").Select()
